$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 676 entirely (post about "「奮闘なくして進歩なし」") - all rows below
# shift up by one as a result.
$ws.Rows.Item(676).Delete()
